$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 31; existing rows 31-86 shift down to 32-87.
$ws.Rows.Item(31).Insert()

# Populate the newly inserted row 31 with the new record.
$ws.Range("A31").Value = 1
$ws.Range("B31").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C31").Value = "Arica y Parinacota"
$ws.Range("D31").Value = 44629
$ws.Range("E31").Value = 15
$ws.Range("F31").Value = "Fruta"
$ws.Range("G31").Value = 100109
$ws.Range("H31").Value = "Uva"
$ws.Range("I31").Value = 100109001
$ws.Range("J31").Value = "Uva"
$ws.Range("K31").Value = "Rosada pastilla"
$ws.Range("L31").Value = "Primera"
$ws.Range("M31").Value = 300
$ws.Range("N31").Value = 14000
$ws.Range("O31").Value = 15000
$ws.Range("P31").Value = 14500
$ws.Range("Q31").Value = "$/bandeja 12 kilos"
$ws.Range("R31").Value = "Región de Coquimbo"
$ws.Range("S31").Value = 1208
$ws.Range("T31").Value = 12
